# Add sign out button in auth view
#
# The "Translation" sheet's Table8 (B3:I799, Text ID / Typography Name /
# Alignment / GB / Direction) loses the two SingleUseId3/SingleUseId4 "UID:
# <value>" / "None" rows plus the two SingleUseId6/SingleUseId7 "Choose your
# avatar" / "essa" rows (original rows 5-8), shifting every following row up
# by four. A brand-new row is then appended at the end of the used data
# (the row that lands on what is now row 14) for a new "Sign out" text
# entry (SingleUseId17 / Default / Center / "Sign out" / LTR).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Remove the four rows (SingleUseId3, SingleUseId4, SingleUseId6, SingleUseId7)
# — this shifts all subsequent rows up by four, same as selecting rows 5:8
# and choosing Delete in Excel.
$ws.Range("A5:A8").EntireRow.Delete() | Out-Null

# Append the new "Sign out" row right after the last remaining data row.
$ws.Cells.Item(14, 2).Value = "SingleUseId17"
$ws.Cells.Item(14, 3).Value = "Default"
$ws.Cells.Item(14, 4).Value = "Center"
$ws.Cells.Item(14, 5).Value = "Sign out"
$ws.Cells.Item(14, 6).Value = "LTR"
